$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1455
$ws.Range("I2").Value = 1542.2858
$ws.Range("K2").Value = 1542.2858
$ws.Range("M2").Value = -1429.2858

$ws.Range("H17").Value = 2690483.8
$ws.Range("J17").Value = 2690483.8
$ws.Range("L17").Value = 8071451.399999999
$ws.Range("N17").Value = -8071787.399999999

$ws.Range("H74").Value = 5498.6875
$ws.Range("I74").Value = 4663.3335
$ws.Range("J74").Value = 5999.9
$ws.Range("K74").Value = 4663.3335
$ws.Range("L74").Value = 5999.9
$ws.Range("M74").Value = -3727.3335
$ws.Range("N74").Value = -7871.9

$ws.Range("H75").Value = 113788.75
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 113788.75
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 113788.75
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -115660.75

$ws.Range("H76").Value = 11006.75
$ws.Range("I76").Value = 11150
$ws.Range("J76").Value = 10004
$ws.Range("K76").Value = 11150
$ws.Range("L76").Value = 10004
$ws.Range("M76").Value = -10835
$ws.Range("N76").Value = -10634

$ws.Range("H77").Value = 5498.6875
$ws.Range("I77").Value = 4663.3335
$ws.Range("J77").Value = 5999.9
$ws.Range("K77").Value = 23316.6675
$ws.Range("L77").Value = 29999.5
$ws.Range("M77").Value = -18636.6675
$ws.Range("N77").Value = -39359.5

$ws.Range("H78").Value = 113788.75
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 113788.75
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 341366.25
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -350726.25

$ws.Range("H79").Value = 11006.75
$ws.Range("I79").Value = 11150
$ws.Range("J79").Value = 10004
$ws.Range("K79").Value = 11150
$ws.Range("L79").Value = 10004
$ws.Range("M79").Value = -10058
$ws.Range("N79").Value = -12188

$ws.Range("H100").Value = 2410.9
$ws.Range("I100").Value = 2157.7144
$ws.Range("K100").Value = 2157.7144
$ws.Range("M100").Value = -1616.7144

$ws.Range("H137").Value = 15499.827
$ws.Range("I137").Value = 14560.4
$ws.Range("J137").Value = 21371.25
$ws.Range("K137").Value = 43681.2
$ws.Range("L137").Value = 64113.75
$ws.Range("M137").Value = -41131.2
$ws.Range("N137").Value = -69213.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4791.06
$ws.Range("I32").Value = 4777.8296
$ws.Range("K32").Value = 4777.8296
$ws.Range("M32").Value = -4490.8296

$ws.Range("H110").Value = 1269.091
$ws.Range("I110").Value = 1066.125
$ws.Range("J110").Value = 1810.3334
$ws.Range("K110").Value = 1066.125
$ws.Range("L110").Value = 1810.3334
$ws.Range("M110").Value = 978.875
$ws.Range("N110").Value = -5900.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3681.6792
$ws.Range("I31").Value = 2326.1428
$ws.Range("K31").Value = 2326.1428
$ws.Range("M31").Value = -2031.1428

$ws.Range("H34").Value = 3681.6792
$ws.Range("I34").Value = 2326.1428
$ws.Range("K34").Value = 2326.1428
$ws.Range("M34").Value = -2124.1428

$ws.Range("H86").Value = 11997.6
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 11997.6
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 11997.6
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -14243.6

$ws.Range("H89").Value = 11997.6
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 11997.6
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 59988
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -71220

$ws.Range("H99").Value = 7500
$ws.Range("I99").Value = 7500
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 7500
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -6002
$ws.Range("N99").ClearContents()

$ws.Range("H106").Value = 15832.667
$ws.Range("J106").Value = 15832.667
$ws.Range("L106").Value = 15832.667
$ws.Range("N106").Value = -18356.667

$ws.Range("H126").Value = 7500
$ws.Range("I126").Value = 7500
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 22500
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -20030
$ws.Range("N126").ClearContents()

$ws.Range("H131").Value = 27615.125
$ws.Range("J131").Value = 30131.572
$ws.Range("L131").Value = 30131.572
$ws.Range("N131").Value = -40211.572

$ws.Range("H132").Value = 3469.6938
$ws.Range("I132").Value = 3281.4285
$ws.Range("K132").Value = 9844.2855
$ws.Range("M132").Value = -7314.2855

$ws.Range("H134").Value = 2520.3635
$ws.Range("I134").Value = 1838.6428
$ws.Range("K134").Value = 5515.928400000001
$ws.Range("M134").Value = -2980.928400000001

$ws.Range("H137").Value = 142334.75
$ws.Range("J137").Value = 149779.67
$ws.Range("L137").Value = 149779.67
$ws.Range("N137").Value = -159979.67

$ws.Range("H141").Value = 85395.5
$ws.Range("J141").Value = 90995
$ws.Range("L141").Value = 90995
$ws.Range("N141").Value = -101355

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 22224220
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 22224220
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 66672660
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -66672916

$ws.Range("H69").Value = 1000
$ws.Range("I69").Value = 1000
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 3000
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -2189
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 1000
$ws.Range("I72").Value = 1000
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 9000
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -4944
$ws.Range("N72").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1263
$ws.Range("I97").Value = 1110.1818
$ws.Range("K97").Value = 1110.1818
$ws.Range("M97").Value = -614.1818000000001

$ws.Range("H132").Value = 5000
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7583.3335
$ws.Range("I7").Value = 6498.3335
$ws.Range("K7").Value = 6498.3335
$ws.Range("M7").Value = -6386.3335

$ws.Range("H40").Value = 4502.5454
$ws.Range("I40").Value = 3902.2
$ws.Range("J40").Value = 5002.8335
$ws.Range("K40").Value = 3902.2
$ws.Range("L40").Value = 5002.8335
$ws.Range("M40").Value = -3766.2
$ws.Range("N40").Value = -5274.8335

$ws.Range("H46").Value = 1971.8148
$ws.Range("I46").Value = 1024.45
$ws.Range("J46").Value = 4678.5713
$ws.Range("K46").Value = 1024.45
$ws.Range("L46").Value = 4678.5713
$ws.Range("M46").Value = -836.45
$ws.Range("N46").Value = -5054.5713

$ws.Range("H61").Value = 4356.8887
$ws.Range("I61").Value = 3901.5
$ws.Range("J61").Value = 8000
$ws.Range("K61").Value = 3901.5
$ws.Range("L61").Value = 8000
$ws.Range("M61").Value = -3699.5
$ws.Range("N61").Value = -8404

$ws.Range("H68").Value = 2966.3333
$ws.Range("I68").Value = 2699
$ws.Range("K68").Value = 2699
$ws.Range("M68").Value = -1950

$ws.Range("H71").Value = 2966.3333
$ws.Range("I71").Value = 2699
$ws.Range("K71").Value = 13495
$ws.Range("M71").Value = -9751

$ws.Range("H113").Value = 4356.8887
$ws.Range("I113").Value = 3901.5
$ws.Range("J113").Value = 8000
$ws.Range("K113").Value = 3901.5
$ws.Range("L113").Value = 8000
$ws.Range("M113").Value = -1731.5
$ws.Range("N113").Value = -12340

$ws.Range("H122").Value = 4409.8
$ws.Range("I122").Value = 4242.7144
$ws.Range("J122").Value = 4799.6665
$ws.Range("K122").Value = 12728.1432
$ws.Range("L122").Value = 14398.9995
$ws.Range("M122").Value = -10278.1432
$ws.Range("N122").Value = -19298.9995

$ws.Range("H126").Value = 7583.3335
$ws.Range("I126").Value = 6498.3335
$ws.Range("K126").Value = 19495.0005
$ws.Range("M126").Value = -17025.0005

$ws.Range("H136").Value = 4050.6924
$ws.Range("I136").Value = 4246.4
$ws.Range("J136").Value = 3398.3333
$ws.Range("K136").Value = 12739.2
$ws.Range("L136").Value = 10194.9999
$ws.Range("M136").Value = -10189.2
$ws.Range("N136").Value = -15294.9999
